$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '27.161.10'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = "'" + '1.899.13'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = "'" + '1.003'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'" + '307.09'
$ws.Range('D6').Value = "'" + '1.003'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = "'" + '0.5241'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = "'" + '0.3806'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('D9').Value = "'" + '0.07310'
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('D10').Value = "'" + '21.33'
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').Value = "'" + '0.9048'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = "'" + '0.08190'
$ws.Range('E12').Value = '  -2.82%  '
$ws.Range('D13').Value = "'" + '1.888.00'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = "'" + '95.26'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = "'" + '5.340'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = "'" + '0.000008665'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = "'" + '14.65'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = "'" + '1.003'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = "'" + '27.192.74'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = "'" + '5.100'
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').Value = "'" + '2.117.24'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').Value = "'" + '10.77'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('D24').Value = "'" + '6.472'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').Value = "'" + '2.333'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'" + '149.63'
$ws.Range('E26').Value = '  +2.02%  '
$ws.Range('D27').Value = "'" + '18.24'
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('D28').Value = "'" + '1.739'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').Value = "'" + '115.23'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').Value = "'" + '4.814'
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('D31').Value = "'" + '4.862'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').Value = "'" + '0.09263'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D33').Value = "'" + '0.05037'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').Value = "'" + '0.7935'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('D35').Value = "'" + '1.218'
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').Value = "'" + '2.966'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = "'" + '3.406'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('D38').Value = "'" + '2.638'
$ws.Range('E38').Value = '  +3.34%  '
$ws.Range('D39').Value = "'" + '0.5712'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').Value = "'" + '0.01995'
$ws.Range('E40').Value = '  +0.99%  '
$ws.Range('D41').Value = "'" + '1.081'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = "'" + '8.989'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'" + '6.600'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').Value = "'" + '116.42'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').Value = "'" + '0.1515'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').Value = "'" + '0.4883'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').Value = "'" + '1.003'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').Value = "'" + '10.13'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').Value = "'" + '1.639'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('D50').Value = "'" + '38.48'
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').Value = "'" + '63.97'
$ws.Range('E51').Value = '  +0.65%  '
